# Applies the "KNX Klemmen Polaritaet vertauscht" edit:
#   1. Merge the two runs "Jung " + "2138" (same formatting) into a single
#      run "Jung 2138" in the Schaltaktor label text box.
#   2. Swap the positions of the "-" and "+" KNX terminal labels.
#   3. (Best effort) refresh the cached NotesMaster date field text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Helper: PowerPoint's Shape.Left/Top are stored as single-precision floats
# (points) and get truncated (not rounded) when converted back to EMU on
# save. Search for a points value whose float32 round-trip truncates to the
# exact EMU we want, so the saved XML offsets match precisely.
# ---------------------------------------------------------------------------
function EmuForPoints($val) {
    $f32 = [System.Single]$val
    $emuF = [double]$f32 * 12700.0
    return [int64]$emuF
}

function FindPointsForEmu($targetEmu) {
    $base = [double]$targetEmu / 12700.0
    if ((EmuForPoints $base) -eq $targetEmu) {
        return $base
    }
    for ($k = 1; $k -lt 5000; $k++) {
        $cand = $base + ($k * 0.0000001)
        if ((EmuForPoints $cand) -eq $targetEmu) {
            return $cand
        }
        $cand2 = $base - ($k * 0.0000001)
        if ((EmuForPoints $cand2) -eq $targetEmu) {
            return $cand2
        }
    }
    return $base
}

function SetShapeOffsetEmu($shape, $xEmu, $yEmu) {
    $shape.Left = FindPointsForEmu $xEmu
    $shape.Top = FindPointsForEmu $yEmu
}

# ---------------------------------------------------------------------------
# 1. Merge "Jung " + "2138" runs into a single run on the Schaltaktor label.
# ---------------------------------------------------------------------------
for ($shpIdx = 1; $shpIdx -le $s.Shapes.Count; $shpIdx++) {
    $sh = $s.Shapes.Item($shpIdx)
    if ($sh.Name -eq "Textfeld 3369991") {
        $tr = $sh.TextFrame.TextRange
        for ($paraIdx = 1; $paraIdx -le $tr.Paragraphs().Count; $paraIdx++) {
            $para = $tr.Paragraphs($paraIdx)
            if ($para.Text -eq "Jung 2138") {
                $run1 = $para.Runs(1)
                # Emptying the first run collapses the run split; the
                # remaining run then gets the full merged text, dropping the
                # now-redundant second run.
                $run1.Text = ""
                $mergedRun = $para.Runs(1)
                $mergedRun.Text = "Jung 2138"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Swap the "-" / "+" KNX terminal label positions (polarity swapped).
# ---------------------------------------------------------------------------
for ($shpIdx2 = 1; $shpIdx2 -le $s.Shapes.Count; $shpIdx2++) {
    $sh = $s.Shapes.Item($shpIdx2)
    if ($sh.Name -eq "Textfeld 171") {
        # was at x=1993072 y=77025
        SetShapeOffsetEmu $sh 2085054 75400
    }
    if ($sh.Name -eq "Textfeld 172") {
        # was at x=2077946 y=84170
        SetShapeOffsetEmu $sh 1986277 88587
    }
}

# ---------------------------------------------------------------------------
# 3. Best effort: refresh cached NotesMaster date field text (read-only in
#    this runtime's object model for master shapes; wrapped so a failure
#    here doesn't abort the rest of the edit).
# ---------------------------------------------------------------------------
try {
    $nm = $p.NotesMaster
    for ($nmIdx = 1; $nmIdx -le $nm.Shapes.Count; $nmIdx++) {
        $sh = $nm.Shapes.Item($nmIdx)
        if ($sh.Name -eq "Datumsplatzhalter 2") {
            $sh.TextFrame.TextRange.Text = "31.10.2020"
        }
    }
} catch {
    # Ignore - not supported by this runtime.
}
